$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("summary")
$wsSummary.Range("C2").Value = 706
$wsSummary.Range("D2").Value = 652
$wsSummary.Range("E2").Value = 85.58
$wsSummary.Range("H2").Value = 1.03
$wsSummary.Range("I2").Value = 0.38
$wsSummary.Range("J2").Value = 0.3
$wsSummary.Range("K2").Value = 0.06
$wsSummary.Range("L2").Value = 1.11
$wsSummary.Range("C3").Value = 706
$wsSummary.Range("D3").Value = 648
$wsSummary.Range("H3").Value = 0.94
$wsSummary.Range("I3").Value = -0.62
$wsSummary.Range("J3").Value = 0.31
$wsSummary.Range("L3").Value = 1.18
$wsSummary.Range("C4").Value = 706
$wsSummary.Range("D4").Value = 598
$wsSummary.Range("E4").Value = 58.7
$wsSummary.Range("F4").Value = -0.43
$wsSummary.Range("G4").Value = 0.1
$wsSummary.Range("I4").Value = 0.45
$wsSummary.Range("J4").Value = 0.38
$wsSummary.Range("K4").Value = 0.05
$wsSummary.Range("C5").Value = 706
$wsSummary.Range("D5").Value = 524
$wsSummary.Range("F5").Value = -0.48
$wsSummary.Range("H5").Value = 0.97
$wsSummary.Range("I5").Value = -0.54
$wsSummary.Range("J5").Value = 0.38
$wsSummary.Range("L5").Value = 0.78
$wsSummary.Range("C6").Value = 706
$wsSummary.Range("D6").Value = 361
$wsSummary.Range("E6").Value = 24.93
$wsSummary.Range("F6").Value = 1.5
$wsSummary.Range("H6").Value = 1.09
$wsSummary.Range("I6").Value = 1.3
$wsSummary.Range("J6").Value = 0.25
$wsSummary.Range("L6").Value = 0.83
$wsSummary.Range("C7").Value = 706
$wsSummary.Range("D7").Value = 671
$wsSummary.Range("E7").Value = 39.64
$wsSummary.Range("F7").Value = 0.56
$wsSummary.Range("H7").Value = 1.04
$wsSummary.Range("I7").Value = 0.99
$wsSummary.Range("J7").Value = 0.32
$wsSummary.Range("L7").Value = 1.06
$wsSummary.Range("C8").Value = 706
$wsSummary.Range("D8").Value = 678
$wsSummary.Range("E8").Value = 37.46
$wsSummary.Range("F8").Value = 0.69
$wsSummary.Range("H8").Value = 1.07
$wsSummary.Range("I8").Value = 1.67
$wsSummary.Range("J8").Value = 0.28
$wsSummary.Range("K8").Value = 0.08
$wsSummary.Range("L8").Value = 0.93
$wsSummary.Range("C9").Value = 706
$wsSummary.Range("D9").Value = 678
$wsSummary.Range("I9").Value = -1.3
$wsSummary.Range("J9").Value = 0.36
$wsSummary.Range("L9").Value = 0.85
$wsSummary.Range("C10").Value = 706
$wsSummary.Range("D10").Value = 676
$wsSummary.Range("E10").Value = 26.48
$wsSummary.Range("F10").Value = 1.33
$wsSummary.Range("H10").Value = 0.95
$wsSummary.Range("I10").Value = -0.96
$wsSummary.Range("K10").Value = 0.05
$wsSummary.Range("L10").Value = 1.61
$wsSummary.Range("C11").Value = 706
$wsSummary.Range("D11").Value = 676
$wsSummary.Range("E11").Value = 23.08
$wsSummary.Range("I11").Value = -0.23
$wsSummary.Range("L11").Value = 1.35

$wsModelFit = $wb.Worksheets.Item("model_fit")
$wsModelFit.Range("B2").Value = 706
$wsModelFit.Range("D2").Value = 7690
$wsModelFit.Range("E2").Value = 7718
$wsModelFit.Range("F2").Value = 7782
$wsModelFit.Range("B3").Value = 706
$wsModelFit.Range("D3").Value = 7655
$wsModelFit.Range("E3").Value = 7701
$wsModelFit.Range("F3").Value = 7806
$wsModelFit.Range("H3").Value = 0.611

$wsSteps = $wb.Worksheets.Item("steps")
$wsSteps.Range("B3").Value = "0.48 (0.112)"
$wsSteps.Range("C3").NumberFormat = "@"
$wsSteps.Range("C3").Value = "-0.48"
